{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the copyright/footer line that followed it, and the blank paragraph that\n// separated them from the \"Requisitos\" list above.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items,text\");\nawait context.sync();\n\nconst items = paras.items;\n\n// Locate the footer paragraphs by their text so the edit is resilient to\n// any paragraphs that may already have been removed/added elsewhere.\nconst jupiterIdx = items.findIndex(p => p.text === \"Ver no Jupiter Salvar em pdf Salvar em docx\");\nconst copyrightIdx = items.findIndex(p => p.text.indexOf(\"Powered by Jekyll and Github pages\") !== -1);\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // The blank paragraph immediately preceding the \"Ver no Jupiter...\" line\n  // is also part of the removed block.\n  const blankIdx = jupiterIdx - 1;\n\n  if (blankIdx >= 0 && items[blankIdx].text === \"\") {\n    items[blankIdx].delete();\n  }\n  items[jupiterIdx].delete();\n  items[copyrightIdx].delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the copyright/footer line that follows it, and the blank paragraph that\n# separates them from the \"Requisitos\" list above.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$targets = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Ver no Jupiter*\" -or $t -like \"*Powered by Jekyll*\") {\n        $targets += $i\n    }\n}\n\nif ($targets.Count -gt 0) {\n    # The blank paragraph immediately preceding the first matched paragraph\n    # is also part of the removed block.\n    $first = $targets[0]\n    if ($first -gt 1) {\n        $prevText = $d.Paragraphs.Item($first - 1).Range.Text\n        if ($prevText -eq \"`r\") {\n            $targets = @($first - 1) + $targets\n        }\n    }\n\n    # Delete from the highest index down so earlier indices stay valid\n    # while the later ones are removed.\n    $sorted = $targets | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
